# Updating single register template
# Add a new data row (row 5) to the "testreg1" worksheet, mirroring the
# existing rows (2-4) of register entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testreg1")

$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "d"
$ws.Cells.Item(5, 3).Value = "d"
$ws.Cells.Item(5, 4).Value = "d"
$ws.Cells.Item(5, 5).Value = "d"
$ws.Cells.Item(5, 6).Value = "Blue book"
$ws.Cells.Item(5, 7).Value = 3

$ws.Range("G6").Select()
